$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new "Time" values (pixel-quantized; closest achievable to ~19.09)
$ws.Columns("A").ColumnWidth = 18.4

# Row 5 previously held text data in B5:G5 - remove it entirely, keep only a date formula in A5
$ws.Range("B5:G5").Clear()
$ws.Range("A5").Formula = "=A4 + 1"

# Row 6 previously held text data in B6:G6 - it gets replaced below with new numeric data
$ws.Range("A6").Formula = "=A5 + 94"

# Rows 7-12 are brand new rows continuing the date series
$ws.Range("A7:A12").Formula = "=A6 + 1"

# Build the new (fill-less, centered, wrapped) numeric style once on a scratch cell, then
# stamp it across the whole new data block B6:G12 in a single paste so every cell shares
# one resolved style entry instead of generating one per assignment.
$ws.Range("Z1").ClearFormats()
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").Copy()
$ws.Range("B6:G12").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Range("B6").Value = 63.4
$ws.Range("C6").Value = 59.9
$ws.Range("D6").Value = 89.2
$ws.Range("E6").Value = 12.7
$ws.Range("F6").Value = 30.1
$ws.Range("G6").Value = 0

$ws.Range("B7").Value = 66.8
$ws.Range("C7").Value = 60.2
$ws.Range("D7").Value = 80.900000000000006
$ws.Range("E7").Value = 8.4
$ws.Range("F7").Value = 30.1
$ws.Range("G7").Value = 0

$ws.Range("B8").Value = 68.599999999999994
$ws.Range("C8").Value = 58.8
$ws.Range("D8").Value = 71.8
$ws.Range("E8").Value = 7.2
$ws.Range("F8").Value = 30.1
$ws.Range("G8").Value = 0

$ws.Range("B9").Value = 67.8
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 76.599999999999994
$ws.Range("E9").Value = 9.5
$ws.Range("F9").Value = 30
$ws.Range("G9").Value = 0

$ws.Range("B10").Value = 67.8
$ws.Range("C10").Value = 59
$ws.Range("D10").Value = 74.3
$ws.Range("E10").Value = 6.5
$ws.Range("F10").Value = 29.9
$ws.Range("G10").Value = 0

$ws.Range("B11").Value = 66.400000000000006
$ws.Range("C11").Value = 59.1
$ws.Range("D11").Value = 78.3
$ws.Range("E11").Value = 9.8000000000000007
$ws.Range("F11").Value = 29.6
$ws.Range("G11").Value = 0

$ws.Range("B12").Value = 67.5
$ws.Range("C12").Value = 49.4
$ws.Range("D12").Value = 54
$ws.Range("E12").Value = 19.8
$ws.Range("F12").Value = 29.5
$ws.Range("G12").Value = 0

# Selection moves to A7, matching the saved sheet view state
$ws.Range("A7").Select()
